# Apply the update to the "Krankenversicherungsbeitraege" workbook:
# - Replace "Beitragsbemessungsgrenze GKV Ost" / "Beitragsbemessungsgrenze GKV West"
#   rows with a single "Beitragsbemessungsgrenze GKV" row and a new
#   "Jahresarbeitsentgeltgrenze GKV" row.
# - Keep the values (74563.82 / 77234.21) as-is, just relabeled.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Tabelle1 - rewrite labels in column A (values in column B stay the same)
$ws1.Range("A2").Value = "ermäßigter Beitragssatz"
$ws1.Range("A3").Value = "Arbeitgeberbeitrag gesetzliche Krankenversicherung in Prozent"
$ws1.Range("A4").Value = "Arbeitnehmerbeitrag gesetzliche Krankenversicherung in Prozent"
$ws1.Range("A5").Value = "Beitragsbemessungsgrenze GKV"
$ws1.Range("A6").Value = "Jahresarbeitsentgeltgrenze GKV"
$ws1.Range("A7").Value = "gueltig_ab"

# Update selection to A9 (was B9)
$ws1.Range("A9").Select()
